$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of price-tracking data to append below the existing table (rows 69-76)
$rows = @(
    @(45215, "19:29", 2573.32, "amazon",        "preto"),
    @(45215, "19:31", 2663,    "mercado livre", "preto"),
    @(45217, "19:33", 2576,    "amazon",        "preto"),
    @(45217, "19:33", 2899,    "mercado livre", "preto"),
    @(45218, "21:24", 2600,    "amazon",        "preto"),
    @(45218, "21:25", 2899,    "mercado livre", "preto"),
    @(45220, "10:26", 3384,    "amazon",        "preto"),
    @(45220, "10:26", 2899,    "mercado livre", "preto")
)

$startRow = 69
$dateFormat = $ws.Cells.Item($startRow - 1, 1).NumberFormat

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
